# Handled stripe kyc error and changed waits while importing products
#
# Replace the old "sanity*" test e-mail addresses used as user credentials
# with the new "*automation*" addresses on each of the credential sheets,
# and move the sheet-5 ("List of Products") selection forward one row.

$wb = $excel.ActiveWorkbook

# Network sheet
$wsNetwork = $wb.Worksheets.Item("Network")
$wsNetwork.Range("A2").Value = "networkautomation1@mailinator.com"

# Vendor 1 sheet
$wsVendor1 = $wb.Worksheets.Item("Vendor 1")
$wsVendor1.Range("A2").Value = "vendorautomation@mailinator.com"

# Vendor 2 sheet
$wsVendor2 = $wb.Worksheets.Item("Vendor 2")
$wsVendor2.Range("A2").Value = "vendorautomation2@mailinator.com"

# Coseller sheet
$wsCoseller = $wb.Worksheets.Item("Coseller")
$wsCoseller.Range("A2").Value = "cosellerautomation@mailinator.com"

# List of Products sheet - move the active selection from A5 to A6
$wsProducts = $wb.Worksheets.Item("List of Products")
$wsProducts.Activate()
$wsProducts.Range("A6").Select()
